# Update the "想去人数" (want-to-go count) figures that changed between
# the two published snapshots of the 丽水-漫展信息 workbook.
#
# Both the "展览" sheet and the "全部类型" sheet carry the same two rows
# of data, so the same F2/F3 update applies to each:
#   F2: 468 -> 469
#   F3: 51  -> 53

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 469
    $ws.Range("F3").Value = 53
}
